# This workbook tracks weekly price-report rows (Brocoli, Femacal de La Calera).
# A new weekly pair of rows (Primera / Segunda) is inserted right after row 958,
# which shifts every later row down by two positions. The two rows that fall off
# the bottom of the original range are appended as brand-new rows 1037/1038.
#
# Strategy: walk the source rows from the bottom (1036) up to the first row that
# needs to move (959), copying each one two rows further down. Because we start
# at the highest row number and move downward, every source row is read before
# it is ever overwritten. Afterwards, rows 959/960 are populated with the new
# week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstShiftRow = 959
$lastRow = 1036
$shiftBy = 2
$dateCol = 4

for ($s = $lastRow; $s -ge $firstShiftRow; $s--) {
    $d = $s + $shiftBy
    for ($c = 1; $c -le 18; $c++) {
        $ws.Cells.Item($d, $c).Value2 = $ws.Cells.Item($s, $c).Value2
    }
    # Column D carries a date-ish number format; keep it consistent on the new row.
    $ws.Cells.Item($d, $dateCol).NumberFormat = $ws.Cells.Item($s, $dateCol).NumberFormat
}

# New week's data, inserted at rows 959 (Primera) and 960 (Segunda).
$ws.Cells.Item(959, 1).Value2 = 3
$ws.Cells.Item(959, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(959, 3).Value2 = "Coquimbo"
$ws.Cells.Item(959, 4).Value2 = 45013
$ws.Cells.Item(959, 5).Value2 = 5
$ws.Cells.Item(959, 6).Value2 = 100112023
$ws.Cells.Item(959, 7).Value2 = "Brócoli"
$ws.Cells.Item(959, 8).Value2 = "Sin especificar"
$ws.Cells.Item(959, 9).Value2 = "Primera"
$ws.Cells.Item(959, 10).Value2 = 3100
$ws.Cells.Item(959, 11).Value2 = 950
$ws.Cells.Item(959, 12).Value2 = 1000
$ws.Cells.Item(959, 13).Value2 = 974
$ws.Cells.Item(959, 14).Value2 = "`$/unidad"
$ws.Cells.Item(959, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(959, 16).Value2 = 974
$ws.Cells.Item(959, 17).Value2 = 1
$ws.Cells.Item(959, 18).Value2 = "Hortaliza"

$ws.Cells.Item(960, 1).Value2 = 3
$ws.Cells.Item(960, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(960, 3).Value2 = "Coquimbo"
$ws.Cells.Item(960, 4).Value2 = 45013
$ws.Cells.Item(960, 5).Value2 = 5
$ws.Cells.Item(960, 6).Value2 = 100112023
$ws.Cells.Item(960, 7).Value2 = "Brócoli"
$ws.Cells.Item(960, 8).Value2 = "Sin especificar"
$ws.Cells.Item(960, 9).Value2 = "Segunda"
$ws.Cells.Item(960, 10).Value2 = 1200
$ws.Cells.Item(960, 11).Value2 = 850
$ws.Cells.Item(960, 12).Value2 = 850
$ws.Cells.Item(960, 13).Value2 = 850
$ws.Cells.Item(960, 14).Value2 = "`$/unidad"
$ws.Cells.Item(960, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(960, 16).Value2 = 850
$ws.Cells.Item(960, 17).Value2 = 1
$ws.Cells.Item(960, 18).Value2 = "Hortaliza"

# Column D keeps the custom date-ish number format on the new rows too.
$dFormat = $ws.Cells.Item(961, $dateCol).NumberFormat
$ws.Cells.Item(959, $dateCol).NumberFormat = $dFormat
$ws.Cells.Item(960, $dateCol).NumberFormat = $dFormat
